$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Brazilian currency code typo: "BRA" -> "BRL" (row 14),
# keep "MXN" on row 15.
$ws.Range("A14").Value = "BRL"
$ws.Range("A15").Value = "MXN"

# Update the saved selection to A15, as reflected in the sheet view.
$ws.Range("A15").Select()
